# Rename COMP_MEAN_WEIGHT* parameter labels to COMP_WEIGHT_MEAN* on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F36").Value = "COMP_WEIGHT_MEAN"
$ws.Range("F44").Value = "COMP_WEIGHT_MEAN_MIDDLE"
$ws.Range("F47").Value = "COMP_WEIGHT_MEAN_CV_MIDDLE"
$ws.Range("F49").Value = "COMP_WEIGHT_MEAN_END"
$ws.Range("F52").Value = "COMP_WEIGHT_MEAN_CV_END"
$ws.Range("F54").Value = "COMP_WEIGHT_MEAN_START"
$ws.Range("F57").Value = "COMP_WEIGHT_MEAN_CV_START"

# Move the active selection to F1, matching the saved cursor position.
$ws.Activate()
$ws.Range("F1").Select()

Write-Host "Renamed COMP_MEAN_WEIGHT to COMP_WEIGHT_MEAN across 7 cells on Sheet1."
